# 417: Add case details for OM_CODE 1003
# Adds a new "case" row (OM_Key 1003, case ref X555555, tier 3, team WMT,
# grade NPSQ, location Community) to each of the four Flag_* report sheets,
# mirroring the row already present on the Court_Reports / Inst_Reports /
# WMT_Extract sheets for OM_CODE 1003.

$wb = $excel.ActiveWorkbook

# --- Flag_Warr_4_n: new row 4 ---
$ws4 = $wb.Worksheets.Item("Flag_Warr_4_n")
$ws4.Range("A2:G2").Copy($ws4.Range("A4:G4"))
$ws4.Range("A4").Value = "I"
$ws4.Range("B4").Value = "X555555"
$ws4.Range("C4").Value = "3"
$ws4.Range("D4").Value = "WMT"
$ws4.Range("E4").Value = "NPSQ"
$ws4.Range("F4").Value = "1003"
$ws4.Range("G4").Value = "Community"
$ws4.Range("A4:G4").Select()

# --- Flag_Upw: new row 4 ---
$ws5 = $wb.Worksheets.Item("Flag_Upw")
$ws5.Range("A2:G2").Copy($ws5.Range("A4:G4"))
$ws5.Range("A4").Value = "I"
$ws5.Range("B4").Value = "X555555"
$ws5.Range("C4").Value = "3"
$ws5.Range("D4").Value = "WMT"
$ws5.Range("E4").Value = "NPSQ"
$ws5.Range("F4").Value = "1003"
$ws5.Range("G4").Value = "Community"
$ws5.Range("A4:G4").Select()

# --- Flag_O_Due: new row 4 ---
$ws6 = $wb.Worksheets.Item("Flag_O_Due")
$ws6.Range("A2:G2").Copy($ws6.Range("A4:G4"))
$ws6.Range("A4").Value = "I"
$ws6.Range("B4").Value = "X555555"
$ws6.Range("C4").Value = "3"
$ws6.Range("D4").Value = "WMT"
$ws6.Range("E4").Value = "NPSQ"
$ws6.Range("F4").Value = "1003"
$ws6.Range("G4").Value = "Community"
$ws6.Range("F4").Select()

# --- Flag_Priority: new row 4 (data) + row 5 (blank formatted trailer row) ---
$ws7 = $wb.Worksheets.Item("Flag_Priority")
$ws7.Range("A3:G3").Copy($ws7.Range("A5:G5"))
$ws7.Range("A5:G5").ClearContents()
$ws7.Range("A2:G2").Copy($ws7.Range("A4:G4"))
$ws7.Range("A4").Value = "I"
$ws7.Range("B4").Value = "X555555"
$ws7.Range("C4").Value = "3"
$ws7.Range("D4").Value = "WMT"
$ws7.Range("E4").Value = "NPSQ"
$ws7.Range("F4").Value = "1003"
$ws7.Range("G4").Value = "Community"
$ws7.Range("A4:G4").Select()

# WMT_Extract becomes the active tab again (was Flag_Priority before the edit)
$ws1 = $wb.Worksheets.Item("WMT_Extract")
$ws1.Activate()
